# Add a "Tissue PC" test to the httk benchmarks table.
# This inserts two new columns (RMSLE.TissuePC, N.TissuePC) into Table1,
# positioned immediately before the existing "Notes" column, and fills in
# the benchmark values for each version row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# --- 1. Physically insert two blank columns right before the "Notes" column (P) ---
# This shifts the existing "Notes" column (P) two slots to the right (to R),
# carrying its data/styles along, and leaves P:Q blank with the same style.
$ws.Range("P1:Q1").EntireColumn.Insert()

# Match the column width used by the other small numeric columns (K:O).
$ws.Range("P1:Q1").EntireColumn.ColumnWidth = 5.8

# --- 2. Grow the table definition to cover the two new columns ---
$tbl.Resize($ws.Range("A1:R23"))

# --- 3. Header labels for the two new columns (also renames the ListColumns) ---
$ws.Range("P1").Value = "RMSLE.TissuePC"
$ws.Range("Q1").Value = "N.TissuePC"
# Touch the Notes header too so the table metadata re-binds its column name/order.
$ws.Range("R1").Value = "Notes"

# --- 4. Fill in the RMSLE.TissuePC / N.TissuePC benchmark values per row ---
$ws.Range("P4").Value = 0.4612
$ws.Range("Q4").Value = 12

$ws.Range("P5").Value = 0.4612
$ws.Range("Q5").Value = 12

$ws.Range("P6").Value = 0.5563
$ws.Range("Q6").Value = 412

$ws.Range("P7").Value = 0.5925
$ws.Range("Q7").Value = 964

$ws.Range("P8").Value = 0.5926
$ws.Range("Q8").Value = 964

$ws.Range("P9").Value = 0.5925
$ws.Range("Q9").Value = 964

$ws.Range("P10").Value = 0.6136
$ws.Range("Q10").Value = 953

$ws.Range("P11").Value = 0.6136
$ws.Range("Q11").Value = 953

$ws.Range("P12").Value = 0.6136
$ws.Range("Q12").Value = 953

$ws.Range("P13").Value = 0.6115
$ws.Range("Q13").Value = 964

$ws.Range("P14").Value = 0.6115
$ws.Range("Q14").Value = 964

$ws.Range("P15").Value = 0.6098
$ws.Range("Q15").Value = 858

$ws.Range("P16").Value = 0.7611
$ws.Range("Q16").Value = 858

$ws.Range("P17").Value = 0.7611
$ws.Range("Q17").Value = 858

$ws.Range("P18").Value = 0.7854
$ws.Range("Q18").Value = 851

$ws.Range("P19").Value = 0.7866
$ws.Range("Q19").Value = 840

$ws.Range("P20").Value = 0.5995
$ws.Range("Q20").Value = 863

$ws.Range("P21").Value = 0.6428
$ws.Range("Q21").Value = 863

$ws.Range("P22").Value = 0.643
$ws.Range("Q22").Value = 863

$ws.Range("P23").Value = 0.63
$ws.Range("Q23").Value = 863

# --- 5. Restore selection/active cell like the authored edit (clicking in the new area) ---
$ws.Activate()
$ws.Range("R4").Select()
